# Auto-generated Excel COM-interop edit script
# Applies the "Updated capital structure database" diff to the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# --- Row 2 (company "3") ---
$ws.Range("D2").Value2 = 0.10705
$ws.Range("E2").Value2 = 0.04550000000000001
$ws.Range("I2").Value2 = 0
$ws.Range("J2").Value2 = 0
$ws.Range("K2").Value2 = 551.9299999999999
$ws.Range("L2").Value2 = 0.3373655256723716
$ws.Range("M2").Value2 = 294.4575
$ws.Range("N2").Value2 = 0.04998599511102059
$ws.Range("O2").Value2 = 0.5335051546391754
$ws.Range("P2").Value2 = 275.4575
$ws.Range("Q2").Value2 = 0.04676062674000137
$ws.Range("R2").Value2 = 0.4990804993386844
$ws.Range("S2").Value2 = 19
$ws.Range("T2").Value2 = 0.06452544085309424
$ws.Range("U2").Value2 = 2933.6
$ws.Range("V2").Value2 = 0.4979968764853671
$ws.Range("W2").Value2 = 0.1314602467257214
$ws.Range("X2").Value2 = 0.05009939762436894
$ws.Range("Y2").Value2 = 0.08136084910135249
$ws.Range("Z2").Value2 = 0.9448182264445151
$ws.Range("AA2").Value2 = 0
$ws.Range("AB2").Value2 = 0.04615201466935195
$ws.Range("AC2").Value2 = -0.04615201466935195
$ws.Range("AD2").Value2 = 1453.4
$ws.Range("AE2").Value2 = 0
$ws.Range("AF2").Value2 = 1453.4
$ws.Range("AG2").Value2 = -1480.2
$ws.Range("AH2").Value2 = 0.1978976607390866
$ws.Range("AI2").Value2 = 0.2370498434237996
$ws.Range("AJ2").Value2 = -0.3356005985580194
$ws.Range("AK2").Value2 = -0.4629096822616962

# --- Row 3 (Banca Transilvania -> BRD - Groupe Société Générale) ---
$ws.Range("B3").Value2 = "BRD - Groupe Société Générale S.A. (BVB:BRD)"
$ws.Range("D3").Value2 = 0.0861
$ws.Range("E3").Value2 = 0.208
$ws.Range("I3").Value2 = 0
$ws.Range("J3").Value2 = 0
$ws.Range("K3").Value2 = 241.9
$ws.Range("L3").Value2 = 0.3491628175519631
$ws.Range("M3").Value2 = 275.2755
$ws.Range("N3").Value2 = 0.1057897467430153
$ws.Range("O3").Value2 = 1.137972302604382
$ws.Range("P3").Value2 = 275.2755
$ws.Range("Q3").Value2 = 0.1057897467430153
$ws.Range("R3").Value2 = 1.137972302604382
$ws.Range("U3").Value2 = 481.4
$ws.Range("V3").Value2 = 0.185004419507321
$ws.Range("W3").Value2 = 0.1314602467257214
$ws.Range("X3").Value2 = 0.04744530327653679
$ws.Range("Y3").Value2 = 0.08401494344918464
$ws.Range("Z3").Value2 = 0.3889512688075454
$ws.Range("AA3").Value2 = 0
$ws.Range("AB3").Value2 = 0.04503010888708486
$ws.Range("AC3").Value2 = -0.04503010888708486
$ws.Range("AD3").Value2 = 500.3
$ws.Range("AE3").Value2 = 0
$ws.Range("AF3").Value2 = 500.3
$ws.Range("AG3").Value2 = 18.90000000000003
$ws.Range("AH3").Value2 = 0.1612622485817432
$ws.Range("AI3").Value2 = 0.1831728480943141
$ws.Range("AJ3").Value2 = 0.007210988172453275
$ws.Range("AK3").Value2 = 0.008400373349926679

# --- Row 4 (BRD - Groupe Société Générale -> Banca Transilvania) ---
$ws.Range("B4").Value2 = "Banca Transilvania S.A. (BVB:TLV)"
$ws.Range("D4").Value2 = 0.128
$ws.Range("E4").Value2 = -0.117
$ws.Range("I4").Value2 = 0
$ws.Range("J4").Value2 = 0
$ws.Range("K4").Value2 = 309
$ws.Range("L4").Value2 = 0.3401210787011558
$ws.Range("M4").Value2 = 19.182
$ws.Range("N4").Value2 = 0.005964181332006716
$ws.Range("O4").Value2 = 0.06207766990291261
$ws.Range("P4").Value2 = 0.182
$ws.Range("Q4").Value2 = 0.00005658852061438965
$ws.Range("R4").Value2 = 0.0005889967637540453
$ws.Range("S4").Value2 = 19
$ws.Range("T4").Value2 = 0.9905119382754667
$ws.Range("U4").Value2 = 2376.5
$ws.Range("V4").Value2 = 0.7389154903302034
$ws.Range("W4").Value2 = 0.1547320981472208
$ws.Range("X4").Value2 = 0.05009939762436894
$ws.Range("Y4").Value2 = 0.1046327005228519
$ws.Range("Z4").Value2 = -8.489860760676573
$ws.Range("AA4").Value2 = -0
$ws.Range("AB4").Value2 = 0.04615201466935195
$ws.Range("AC4").Value2 = -0.04615201466935195
$ws.Range("AD4").Value2 = 927.8
$ws.Range("AE4").Value2 = 0
$ws.Range("AF4").Value2 = 927.8
$ws.Range("AG4").Value2 = -1448.7
$ws.Range("AH4").Value2 = 0.2238899613899614
$ws.Range("AI4").Value2 = 0.2815098003519631
$ws.Range("AJ4").Value2 = -0.8196322489391797
$ws.Range("AK4").Value2 = -1.575872946807354

# --- Row 5 (Patria Bank SA) ---
$ws.Range("I5").Value2 = 0
$ws.Range("J5").Value2 = 0
$ws.Range("K5").Value2 = 1.03
$ws.Range("L5").Value2 = 0.02968299711815562
$ws.Range("M5").Value2 = -0
$ws.Range("N5").Value2 = -0
$ws.Range("O5").Value2 = -0
$ws.Range("S5").Value2 = 0
$ws.Range("U5").Value2 = 75.7
$ws.Range("V5").Value2 = 1.044137931034483
$ws.Range("W5").Value2 = 0.01358839050131926
$ws.Range("X5").Value2 = 0.05176807211860267
$ws.Range("Y5").Value2 = -0.0381796816172834
$ws.Range("Z5").Value2 = 0.6049511854951186
$ws.Range("AA5").Value2 = 0
$ws.Range("AB5").Value2 = 0.04677544024914178
$ws.Range("AC5").Value2 = -0.04677544024914178
$ws.Range("AD5").Value2 = 25.3
$ws.Range("AE5").Value2 = 0
$ws.Range("AF5").Value2 = 25.3
$ws.Range("AG5").Value2 = -50.40000000000001
$ws.Range("AH5").Value2 = 0.2586912065439673
$ws.Range("AI5").Value2 = 0.2430355427473583
$ws.Range("AJ5").Value2 = -2.28054298642534
$ws.Range("AK5").Value2 = -1.774647887323944

# --- Remove now-obsolete cells (AN/AP debt_ebitda columns, extra T5) -----
$ws.Range("AN2").ClearContents()
$ws.Range("AP2").ClearContents()
$ws.Range("AN3").ClearContents()
$ws.Range("AP3").ClearContents()
$ws.Range("AN4").ClearContents()
$ws.Range("AP4").ClearContents()
$ws.Range("T5").ClearContents()
$ws.Range("AN5").ClearContents()
$ws.Range("AP5").ClearContents()

